$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "ノースウィンドトレーダーズは、ワイドワールドインポーターのためのビールと醸造所の優先サプライヤーです。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Northwind Traders は、Wide World Importers のためのビールとサイダーの優先サプライヤーです。",
    2)

$d.Content.Find.Execute(
    "2023 年 2 月 1 日に Wide World Importers と交渉されたサプライヤー契約条件は次のとおりです。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2023 年 2 月 1 日に Wide World Importers と交渉したサプライヤー契約条件は次のとおりです。",
    2)

$d.Content.Find.Execute(
    "10 日以内支払割引 2%、支払期限 30 日以内",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "10 日以内支払割引 2%、支払期限 45 日以内",
    2)

$d.Content.Find.Execute(
    "一律料金",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "遅延料金",
    2)

$d.Content.Find.Execute(
    "1 か月あたり `$100",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1 か月あたり 2%",
    2)

$d.Content.Find.Execute(
    "最小注文金額",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "最小注文量",
    2)

$d.Content.Find.Execute(
    "20 時間/月",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1 か月あたり 50 ケース",
    2)

$d.Content.Find.Execute(
    "最大注文金額",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "最大注文量",
    2)

$d.Content.Find.Execute(
    "0 (最大値なし)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "最大値なし",
    2)

$d.Content.Find.Execute(
    "ケースあたりの価格は `$25 で固定されています",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ケースあたりの価格は 25 ドルで固定されています",
    2)
